# Automatische test-sync: 2025-08-04 20:38:50
#
# Appends the new "Testmail #6" follow-up log entry to the "Logs" sheet
# (row 18) and updates the "Dashboard" summary sheet so the
# "Inkoop / Bestellingen" category count reflects the new row, matching
# the category ordering produced by the upstream export.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new row to the "Logs" sheet -----------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 18
$logs.Cells.Item($newRow, 1).Value = "Hebben we EcoPro-700 nog op voorraad?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #6: Hebben we EcoPro-700 nog op voorraad?"
$logs.Cells.Item($newRow, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-04 20:38:25"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional-formatting ranges (D/G/H/I/J) so they keep
# covering the full data range through the newly added row.
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range($col + "2:" + $col + "17")
    $newRange = $logs.Range($col + "2:" + $col + "18")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count(); $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- 2. Update the "Dashboard" summary sheet -------------------------------
# The new row belongs to "Inkoop / Bestellingen", so its tally grows to 4
# and (matching the authoritative export) it now sorts above
# "Opvolging / Status" in the category list, so the two rows swap places.
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(3, 1).Value = "Inkoop / Bestellingen"
$dash.Cells.Item(3, 2).Value = 4
$dash.Cells.Item(4, 1).Value = "Opvolging / Status"
$dash.Cells.Item(4, 2).Value = 3
